# "Условие лабораторной работы" -> "Задание лабораторной работы"
# (two occurrences, each followed by a period that lives in its own,
# separately-formatted run). A straight text replace collapses
# everything into one run; the committed version instead keeps
# "Задание" / " лабораторной работы" / "." as three sibling runs
# (all sharing the same identical rPr). We do the plain text swap
# first, then re-create the two run boundaries by toggling Bold off
# and back on over just the sub-ranges that must become their own
# run - a no-op formatting-wise, but it forces Word to split the run
# there instead of leaving one big merged run.
$d = $word.ActiveDocument

$d.Content.Find.Execute("Условие лабораторной работы", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Задание лабораторной работы", 2)

$hit = $d.Content
while ($hit.Find.Execute("Задание лабораторной работы.", $true, $false, $false, $false, $false, `
                          $true, 1, $false, $null, 0)) {
    $start = $hit.Start
    $end = $hit.End

    # Boundary between "Задание" and " лабораторной работы"
    $part1 = $d.Range($start, $start + 7)
    $part1.Bold = $false
    $part1.Bold = $true

    # Boundary between " лабораторной работы" and "."
    $part3 = $d.Range($end - 1, $end)
    $part3.Bold = $false
    $part3.Bold = $true

    $hit.Collapse(0)
}
